$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1907216494845361
$ws.Range("C2").Value = 0.5481099656357389
$ws.Range("J2").Value = 0.01202749140893471
$ws.Range("P2").Value = 0.1460481099656357
$ws.Range("S2").Value = 0.1030927835051546
$ws.Range("B3").Value = 0.009202453987730062
$ws.Range("C3").Value = 0.03067484662576687
$ws.Range("J3").Value = 0.03374233128834356
$ws.Range("P3").Value = 0.7177914110429447
$ws.Range("S3").Value = 0.2085889570552147
$ws.Range("J4").Value = 0.03947368421052631
$ws.Range("O4").Value = 0.0131578947368421
$ws.Range("P4").Value = 0.6447368421052632
$ws.Range("S4").Value = 0.3026315789473684
$ws.Range("B6").Value = 0.05825242718446602
$ws.Range("D6").Value = 0.007766990291262136
$ws.Range("E6").Value = 0.001941747572815534
$ws.Range("F6").Value = 0.05048543689320388
$ws.Range("J6").Value = 0.2563106796116505
$ws.Range("O6").Value = 0.01359223300970874
$ws.Range("Q6").Value = 0.1689320388349514
$ws.Range("R6").Value = 0.08737864077669903
$ws.Range("S6").Value = 0.3553398058252427
$ws.Range("B7").Value = 0.1213592233009709
$ws.Range("D7").Value = 0.01456310679611651
$ws.Range("F7").Value = 0.08495145631067962
$ws.Range("J7").Value = 0.1116504854368932
$ws.Range("O7").Value = 0.01941747572815534
$ws.Range("Q7").Value = 0.1966019417475728
$ws.Range("R7").Value = 0.04854368932038835
$ws.Range("S7").Value = 0.4029126213592233
$ws.Range("B8").Value = 0.0872210953346856
$ws.Range("D8").Value = 0.01521298174442191
$ws.Range("E8").Value = 0.002028397565922921
$ws.Range("F8").Value = 0.06490872210953347
$ws.Range("J8").Value = 0.103448275862069
$ws.Range("O8").Value = 0.02738336713995943
$ws.Range("Q8").Value = 0.1977687626774848
$ws.Range("R8").Value = 0.09837728194726167
$ws.Range("S8").Value = 0.4036511156186612
$ws.Range("B9").Value = 0.09510869565217392
$ws.Range("D9").Value = 0.03260869565217391
$ws.Range("E9").Value = 0.002717391304347826
$ws.Range("F9").Value = 0.08152173913043478
$ws.Range("J9").Value = 0.108695652173913
$ws.Range("O9").Value = 0.01902173913043478
$ws.Range("Q9").Value = 0.1902173913043478
$ws.Range("R9").Value = 0.09239130434782608
$ws.Range("S9").Value = 0.3777173913043478
$ws.Range("B10").Value = 0.1074148296593186
$ws.Range("D10").Value = 0.01683366733466934
$ws.Range("E10").Value = 0.0008016032064128256
$ws.Range("F10").Value = 0.08216432865731463
$ws.Range("J10").Value = 0.1042084168336673
$ws.Range("O10").Value = 0.02364729458917836
$ws.Range("Q10").Value = 0.2160320641282565
$ws.Range("R10").Value = 0.08376753507014029
$ws.Range("S10").Value = 0.3651302605210421
$ws.Range("J11").Value = 0.095
$ws.Range("K11").Value = 0.17
$ws.Range("L11").Value = 0.5933333333333334
$ws.Range("S11").Value = 0.01666666666666667
$ws.Range("G12").Value = 0.7331536388140162
$ws.Range("J12").Value = 0.1913746630727763
$ws.Range("K12").Value = 0.005390835579514825
$ws.Range("L12").Value = 0.03773584905660377
$ws.Range("S12").Value = 0.03234501347708895
$ws.Range("G13").Value = 0.6764705882352942
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.02941176470588235
$ws.Range("F15").Value = 0.01803607214428858
$ws.Range("H15").Value = 0.1603206412825651
$ws.Range("I15").Value = 0.04208416833667335
$ws.Range("J15").Value = 0.3827655310621242
$ws.Range("K15").Value = 0.06813627254509018
$ws.Range("M15").Value = 0.01002004008016032
$ws.Range("N15").Value = 0.002004008016032064
$ws.Range("O15").Value = 0.0781563126252505
$ws.Range("S15").Value = 0.2384769539078156
$ws.Range("F16").Value = 0.01971830985915493
$ws.Range("H16").Value = 0.1774647887323944
$ws.Range("I16").Value = 0.07887323943661972
$ws.Range("J16").Value = 0.3943661971830986
$ws.Range("K16").Value = 0.123943661971831
$ws.Range("M16").Value = 0.01690140845070422
$ws.Range("O16").Value = 0.05633802816901409
$ws.Range("S16").Value = 0.1323943661971831
$ws.Range("F17").Value = 0.01859504132231405
$ws.Range("H17").Value = 0.1807851239669422
$ws.Range("I17").Value = 0.08677685950413223
$ws.Range("J17").Value = 0.4297520661157025
$ws.Range("K17").Value = 0.09194214876033058
$ws.Range("M17").Value = 0.02685950413223141
$ws.Range("N17").Value = 0.001033057851239669
$ws.Range("O17").Value = 0.08057851239669421
$ws.Range("S17").Value = 0.08367768595041322
$ws.Range("F18").Value = 0.02227722772277228
$ws.Range("H18").Value = 0.2128712871287129
$ws.Range("I18").Value = 0.09405940594059406
$ws.Range("J18").Value = 0.3514851485148515
$ws.Range("K18").Value = 0.09405940594059406
$ws.Range("M18").Value = 0.0198019801980198
$ws.Range("N18").Value = 0.002475247524752475
$ws.Range("O18").Value = 0.07178217821782178
$ws.Range("S18").Value = 0.1311881188118812
$ws.Range("F19").Value = 0.0168946098149638
$ws.Range("H19").Value = 0.2377312952534192
$ws.Range("I19").Value = 0.08286403861625101
$ws.Range("J19").Value = 0.3527755430410298
$ws.Range("K19").Value = 0.1142397425583266
$ws.Range("M19").Value = 0.02373290426387771
$ws.Range("N19").Value = 0.001609010458567981
$ws.Range("O19").Value = 0.07160096540627514
$ws.Range("S19").Value = 0.09855189058728882
